$d = $word.ActiveDocument
Write-Output $d.Bookmarks.Count
for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
    $b = $d.Bookmarks.Item($i)
    Write-Output ($b.Name + " : " + $b.Start + "-" + $b.End)
}
